$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.443.75'
$ws.Range('D3').Value = '1.946.13'
$ws.Range('E3').Value = '  -0.73%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''243.21'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('D7').Value = '''58.24'
$ws.Range('E7').Value = '  -2.47%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  -2.13%  '
$ws.Range('D10').Value = '''55.84'
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('D11').Value = '''0.0833'
$ws.Range('E11').Value = '  +3.97%  '
$ws.Range('D12').Value = '''0.104'
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('D13').Value = '''0.823'
$ws.Range('E13').Value = '  -3.66%  '
$ws.Range('D14').Value = '''21.50'
$ws.Range('E14').Value = '  -2.61%  '
$ws.Range('D15').Value = '2.230.50'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = '''13.60'
$ws.Range('E16').Value = '  -2.70%  '
$ws.Range('E17').Value = '  -2.79%  '
$ws.Range('D18').Value = '1.964.56'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = '36.353.52'
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').Value = '''69.54'
$ws.Range('E20').Value = '  -1.74%  '
$ws.Range('D21').Value = '0.0₃0860'
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').Value = '''228.62'
$ws.Range('E22').Value = '  -2.62%  '
$ws.Range('D23').Value = '''5.05'
$ws.Range('E23').Value = '  -2.59%  '
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('E25').Value = '  -3.21%  '
$ws.Range('E26').Value = '  +0.69%  '
$ws.Range('E27').Value = '  -5.38%  '
$ws.Range('D28').Value = '''161.63'
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('D29').Value = '''0.131'
$ws.Range('E29').Value = '  +0.90%  '
$ws.Range('D30').Value = '''19.50'
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('E32').Value = '  +1.59%  '
$ws.Range('E33').Value = '  -3.62%  '
$ws.Range('D34').Value = '''0.0629'
$ws.Range('E34').Value = '  +1.98%  '
$ws.Range('E35').Value = '  -2.66%  '
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').Value = '''1.77'
$ws.Range('E38').Value = '  -3.04%  '
$ws.Range('E39').Value = '  -5.84%  '
$ws.Range('D40').Value = '''3.04'
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('D41').Value = '''0.0984'
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('E43').Value = '  -4.09%  '
$ws.Range('E44').Value = '  -1.31%  '
$ws.Range('D45').Value = '''16.00'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('D46').Value = '1.362.05'
$ws.Range('E46').Value = '  +2.30%  '
$ws.Range('E47').Value = '  -4.41%  '
$ws.Range('D48').Value = '''87.97'
$ws.Range('E48').Value = '  -4.15%  '
$ws.Range('E49').Value = '  -4.40%  '
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('D51').Value = '''45.47'
$ws.Range('E51').Value = '  +3.74%  '
